$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw instruction counts (B column) - tracking the data-cache stack
# pointer directly instead of recalculating it.
$ws.Range("B6").Value  = 1019
$ws.Range("B7").Value  = 79
$ws.Range("B8").Value  = 423
$ws.Range("B9").Value  = 277
$ws.Range("B10").Value = 424
$ws.Range("B14").Value = 171
$ws.Range("B15").Value = 198

# Update the view to match where the author was looking when they made the edit.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("H12").Select()
